# Auto-generated Excel COM-interop script applying the Jenova_Profits workbook update.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ /
# LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns (H-N) on affected rows across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51: A Bile Business | Shark Oil
$ws.Range("H51").Value = 2755.2856
$ws.Range("I51").Value = 2957.4
$ws.Range("K51").Value = 2957.4
$ws.Range("M51").Value = -2473.4

# Row 76: Warding Off Temptation | Enchanted Hardsilver Ink
$ws.Range("I76").Value = 7800
$ws.Range("K76").Value = 7800
$ws.Range("M76").Value = -7485

# Row 79: The Garden of Arcane Delights (L) | Enchanted Hardsilver Ink
$ws.Range("I79").Value = 7800
$ws.Range("K79").Value = 7800
$ws.Range("M79").Value = -6708

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 2026.3247
$ws.Range("I132").Value = 1917.942
$ws.Range("J132").Value = 2961.125
$ws.Range("K132").Value = 5753.826
$ws.Range("L132").Value = 8883.375
$ws.Range("M132").Value = -3223.826
$ws.Range("N132").Value = -13943.375

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 10244.409
$ws.Range("J138").Value = 10443.284
$ws.Range("L138").Value = 31329.852
$ws.Range("N138").Value = -41609.852

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 47166.793
$ws.Range("I2").Value = 62243.832
$ws.Range("K2").Value = 62243.832
$ws.Range("M2").Value = -62130.832

# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 22189.3
$ws.Range("I32").Value = 19099.223
$ws.Range("K32").Value = 19099.223
$ws.Range("M32").Value = -18812.223

# Row 43: They've Got Legs | Steel Sabatons
$ws.Range("H43").Value = 19492.445
$ws.Range("J43").Value = 20511.25
$ws.Range("L43").Value = 20511.25
$ws.Range("N43").Value = -21137.25

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 4552.5713
$ws.Range("I74").Value = 3644.6667
$ws.Range("K74").Value = 3644.6667
$ws.Range("M74").Value = -2770.6667

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 4552.5713
$ws.Range("I77").Value = 3644.6667
$ws.Range("K77").Value = 18223.3335
$ws.Range("M77").Value = -13855.3335

# Row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 47166.793
$ws.Range("I116").Value = 62243.832
$ws.Range("K116").Value = 62243.832
$ws.Range("M116").Value = -59949.832

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 4159.1177
$ws.Range("I132").Value = 3389.5144
$ws.Range("K132").Value = 10168.5432
$ws.Range("M132").Value = -7638.5432

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 47166.793
$ws.Range("I3").Value = 62243.832
$ws.Range("K3").Value = 62243.832
$ws.Range("M3").Value = -62129.832

# Row 12: A Hit Job | Bronze Chaser Hammer
$ws.Range("H12").Value = 850
$ws.Range("I12").Value = 200
$ws.Range("K12").Value = 200
$ws.Range("M12").Value = -32

$ws = $wb.Worksheets.Item("CRP")
# Row 4: A Clogful of Camaraderie | Maple Clogs
$ws.Range("H4").Value = 10003813
$ws.Range("I4").Value = 3500
$ws.Range("J4").Value = 16004000
$ws.Range("K4").Value = 3500
$ws.Range("L4").Value = 16004000
$ws.Range("M4").Value = -3388
$ws.Range("N4").Value = -16004224

# Row 13: Compulsory Conjury | Maple Cane
$ws.Range("H13").Value = 4900
$ws.Range("J13").Value = 4900
$ws.Range("L13").Value = 4900
$ws.Range("N13").Value = -5178

# Row 16: Raise the Roof | Ash Lumber
$ws.Range("H16").Value = 21885.777
$ws.Range("I16").Value = 9425.210999999999
$ws.Range("J16").Value = 51479.625
$ws.Range("K16").Value = 9425.210999999999
$ws.Range("L16").Value = 51479.625
$ws.Range("M16").Value = -9138.210999999999
$ws.Range("N16").Value = -52053.625

# Row 19: Shielding Sales | Square Ash Shield
$ws.Range("H19").Value = 1835.55
$ws.Range("I19").Value = 232.2
$ws.Range("J19").Value = 3438.9
$ws.Range("K19").Value = 232.2
$ws.Range("L19").Value = 3438.9
$ws.Range("M19").Value = -62.19999999999999
$ws.Range("N19").Value = -3778.9

# Row 24: What You Need | Square Ash Shield
$ws.Range("H24").Value = 1835.55
$ws.Range("I24").Value = 232.2
$ws.Range("J24").Value = 3438.9
$ws.Range("K24").Value = 232.2
$ws.Range("L24").Value = 3438.9
$ws.Range("M24").Value = -62.19999999999999
$ws.Range("N24").Value = -3778.9

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 66061.17999999999
$ws.Range("I31").Value = 4428.4165
$ws.Range("J31").Value = 213979.8
$ws.Range("K31").Value = 4428.4165
$ws.Range("L31").Value = 213979.8
$ws.Range("M31").Value = -4133.4165
$ws.Range("N31").Value = -214569.8

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 66061.17999999999
$ws.Range("I34").Value = 4428.4165
$ws.Range("J34").Value = 213979.8
$ws.Range("K34").Value = 4428.4165
$ws.Range("L34").Value = 213979.8
$ws.Range("M34").Value = -4226.4165
$ws.Range("N34").Value = -214383.8

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Range("H58").Value = 3121.64
$ws.Range("I58").Value = 3078.5
$ws.Range("J58").Value = 3161.4614
$ws.Range("K58").Value = 3078.5
$ws.Range("L58").Value = 3161.4614
$ws.Range("M58").Value = -2875.5
$ws.Range("N58").Value = -3567.4614

# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 7174.875
$ws.Range("I99").Value = 6800
$ws.Range("K99").Value = 6800
$ws.Range("M99").Value = -5302

# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("H105").Value = 1256
$ws.Range("I105").Value = 781.6
$ws.Range("K105").Value = 781.6
$ws.Range("M105").Value = 965.4

# Row 113: Patient Patients | White Ash Lumber
$ws.Range("H113").Value = 21885.777
$ws.Range("I113").Value = 9425.210999999999
$ws.Range("J113").Value = 51479.625
$ws.Range("K113").Value = 9425.210999999999
$ws.Range("L113").Value = 51479.625
$ws.Range("M113").Value = -7255.210999999999
$ws.Range("N113").Value = -55819.625

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 7174.875
$ws.Range("I126").Value = 6800
$ws.Range("K126").Value = 20400
$ws.Range("M126").Value = -17930

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Range("H136").Value = 3121.64
$ws.Range("I136").Value = 3078.5
$ws.Range("J136").Value = 3161.4614
$ws.Range("K136").Value = 9235.5
$ws.Range("L136").Value = 9484.3842
$ws.Range("M136").Value = -6685.5
$ws.Range("N136").Value = -14584.3842

$ws = $wb.Worksheets.Item("CUL")
# Row 32: Convalescence Precedes Essence | Ginger Cookie
$ws.Range("H32").Value = 11252176
$ws.Range("J32").Value = 15000003
$ws.Range("L32").Value = 45000009
$ws.Range("N32").Value = -45000575

# Row 100: Souper | Gameni
$ws.Range("H100").Value = 2507
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers | Copper Ingot
$ws.Range("H2").Value = 322.5
$ws.Range("I2").Value = 460
$ws.Range("J2").Value = 93.333336
$ws.Range("K2").Value = 460
$ws.Range("L2").Value = 93.333336
$ws.Range("M2").Value = -347
$ws.Range("N2").Value = -319.333336

# Row 11: A Ringing Success | Copper Ring
$ws.Range("H11").Value = 3348444.2
$ws.Range("I11").Value = 5006000
$ws.Range("J11").Value = 33333
$ws.Range("K11").Value = 5006000
$ws.Range("L11").Value = 33333
$ws.Range("M11").Value = -5005861
$ws.Range("N11").Value = -33611

# Row 18: Gorgeous Gorget | Brass Gorget
$ws.Range("H18").Value = 1000000000
$ws.Range("J18").Value = 1000000000
$ws.Range("L18").Value = 1000000000
$ws.Range("N18").Value = -1000000586

# Row 46: Burning the Midnight Oil | Fire Brand
$ws.Range("H46").Value = 41230.46
$ws.Range("J46").Value = 41230.46
$ws.Range("L46").Value = 41230.46
$ws.Range("N46").Value = -41542.46

# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 2860.5518
$ws.Range("I126").Value = 2369.9443
$ws.Range("K126").Value = 7109.8329
$ws.Range("M126").Value = -4639.8329

$ws = $wb.Worksheets.Item("LTW")
# Row 110: Breeches of Trust | Gliderskin Breeches of Fending
$ws.Range("H110").Value = 79734.14
$ws.Range("J110").Value = 79734.14
$ws.Range("L110").Value = 79734.14
$ws.Range("N110").Value = -87914.14

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 5379.909
$ws.Range("I132").Value = 4804.72
$ws.Range("K132").Value = 14414.16
$ws.Range("M132").Value = -11884.16

$ws = $wb.Worksheets.Item("WVR")
# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 1771.5714
$ws.Range("I126").Value = 1550.3
$ws.Range("K126").Value = 4650.9
$ws.Range("M126").Value = -2180.9
